$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "ECs" sending-cluster rows (rows 5-7), keeping only the
# "FAPs" sending-cluster rows, and shrinking the used range to A1:T4.
$ws.Rows("5:7").Delete()

# Refresh data for the remaining rows (2-4) with the updated TPM-derived
# values. Columns A-D hold the category labels, E-T hold the recomputed
# numeric metrics.

# Row 2: FAPs -> Wnt3/Lrp6 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt3"
$ws.Range("C2").Value = "Lrp6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1741663333333333
$ws.Range("H2").Value = 0.5224989999999999
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.40685866666667
$ws.Range("N2").Value = 37.220576
$ws.Range("O2").Value = 0.1720325859617629
$ws.Range("P2").Value = 0.1720325859617629
$ws.Range("Q2").Value = 2.160857082158222
$ws.Range("R2").Value = 19.447713739424
$ws.Range("S2").Value = 0.1720325859617629
$ws.Range("T2").Value = 0.1720325859617629

# Row 3: FAPs -> Wnt3/Lrp6 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt3"
$ws.Range("C3").Value = "Lrp6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1741663333333333
$ws.Range("H3").Value = 0.5224989999999999
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 43.97212233333332
$ws.Range("N3").Value = 131.916367
$ws.Range("O3").Value = 0.6097142007069145
$ws.Range("P3").Value = 0.6097142007069145
$ws.Range("Q3").Value = 7.658463315681442
$ws.Range("R3").Value = 68.92616984113297
$ws.Range("S3").Value = 0.6097142007069145
$ws.Range("T3").Value = 0.6097142007069145

# Row 4: FAPs -> Wnt3/Lrp6 -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt3"
$ws.Range("C4").Value = "Lrp6"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1741663333333333
$ws.Range("H4").Value = 0.5224989999999999
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 15.740255
$ws.Range("N4").Value = 47.220765
$ws.Range("O4").Value = 0.2182532133313226
$ws.Range("P4").Value = 0.2182532133313226
$ws.Range("Q4").Value = 2.741422499081666
$ws.Range("R4").Value = 24.672802491735
$ws.Range("S4").Value = 0.2182532133313226
$ws.Range("T4").Value = 0.2182532133313226
